{"js": "// Replace the two-digit multiplication problems/answers in the table\n// with the new set of values, preserving all run/paragraph formatting.\nconst replacements = [\n  [\"72\u00d776=5472\", \"18\u00d758=1044\"],\n  [\"13\u00d791=1183\", \"37\u00d796=3552\"],\n  [\"67\u00d718=1206\", \"58\u00d719=1102\"],\n  [\"56\u00d799=5544\", \"88\u00d733=2904\"],\n  [\"92\u00d738=3496\", \"93\u00d771=6603\"],\n  [\"93\u00d791=8463\", \"67\u00d713=871\"],\n  [\"85\u00d738=3230\", \"89\u00d795=8455\"],\n  [\"83\u00d779=6557\", \"72\u00d727=1944\"],\n  [\"66\u00d799=6534\", \"16\u00d782=1312\"],\n  [\"67\u00d754=3618\", \"74\u00d779=5846\"],\n  [\"67\u00d769=4623\", \"67\u00d723=1541\"],\n  [\"25\u00d740=1000\", \"18\u00d786=1548\"],\n  [\"37\u00d778=2886\", \"75\u00d732=2400\"],\n  [\"69\u00d795=6555\", \"55\u00d740=2200\"],\n  [\"90\u00d788=7920\", \"70\u00d769=4830\"],\n  [\"12\u00d772=864\", \"17\u00d726=442\"],\n  [\"21\u00d797=2037\", \"35\u00d728=980\"],\n  [\"72\u00d798=7056\", \"63\u00d734=2142\"],\n  [\"73\u00d734=2482\", \"17\u00d792=1564\"],\n  [\"53\u00d716=848\", \"25\u00d719=475\"],\n  [\"41\u00d784=3444\", \"82\u00d797=7954\"],\n  [\"50\u00d788=4400\", \"19\u00d724=456\"],\n  [\"46\u00d716=736\", \"97\u00d779=7663\"],\n  [\"43\u00d745=1935\", \"91\u00d753=4823\"],\n  [\"52\u00d745=2340\", \"25\u00d774=1850\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems/answers in the table\n# with the new set of values, preserving all run/paragraph formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"72\u00d776=5472\", \"18\u00d758=1044\"),\n    @(\"13\u00d791=1183\", \"37\u00d796=3552\"),\n    @(\"67\u00d718=1206\", \"58\u00d719=1102\"),\n    @(\"56\u00d799=5544\", \"88\u00d733=2904\"),\n    @(\"92\u00d738=3496\", \"93\u00d771=6603\"),\n    @(\"93\u00d791=8463\", \"67\u00d713=871\"),\n    @(\"85\u00d738=3230\", \"89\u00d795=8455\"),\n    @(\"83\u00d779=6557\", \"72\u00d727=1944\"),\n    @(\"66\u00d799=6534\", \"16\u00d782=1312\"),\n    @(\"67\u00d754=3618\", \"74\u00d779=5846\"),\n    @(\"67\u00d769=4623\", \"67\u00d723=1541\"),\n    @(\"25\u00d740=1000\", \"18\u00d786=1548\"),\n    @(\"37\u00d778=2886\", \"75\u00d732=2400\"),\n    @(\"69\u00d795=6555\", \"55\u00d740=2200\"),\n    @(\"90\u00d788=7920\", \"70\u00d769=4830\"),\n    @(\"12\u00d772=864\",  \"17\u00d726=442\"),\n    @(\"21\u00d797=2037\", \"35\u00d728=980\"),\n    @(\"72\u00d798=7056\", \"63\u00d734=2142\"),\n    @(\"73\u00d734=2482\", \"17\u00d792=1564\"),\n    @(\"53\u00d716=848\",  \"25\u00d719=475\"),\n    @(\"41\u00d784=3444\", \"82\u00d797=7954\"),\n    @(\"50\u00d788=4400\", \"19\u00d724=456\"),\n    @(\"46\u00d716=736\",  \"97\u00d779=7663\"),\n    @(\"43\u00d745=1935\", \"91\u00d753=4823\"),\n    @(\"52\u00d745=2340\", \"25\u00d774=1850\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
